$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = -7.995
$ws.Range("D13").Value = -8.147
$ws.Range("D16").Value = -8.563999999999998
$ws.Range("D18").Value = -8.27
$ws.Range("D20").Value = -7.691000000000001
$ws.Range("D26").Value = -7.533000000000001
$ws.Range("D27").Value = -8.056000000000001
$ws.Range("D29").Value = -7.382
$ws.Range("D35").Value = -7.888
$ws.Range("D36").Value = -7.540999999999999
$ws.Range("D45").Value = -7.556999999999999
$ws.Range("D55").Value = -8.300000000000001
$ws.Range("D57").Value = -8.27
$ws.Range("D69").Value = -7.406000000000001
$ws.Range("D76").Value = -7.943000000000001
$ws.Range("D78").Value = -8.077000000000002
$ws.Range("D82").Value = -8.17
$ws.Range("D83").Value = -8.103
$ws.Range("D93").Value = -6.976999999999999
$ws.Range("D97").Value = -7.423
